$d = $word.ActiveDocument

# --- Define the three new character styles (added at the end of styles.xml) ---

$gaNStyle = $d.Styles.Add("GaNStyle", 2)            # wdStyleTypeCharacter
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)    # wdStyleTypeCharacter
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)            # wdStyleTypeCharacter
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608                      # BGR for 000080 (navy)
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1                        # wdUnderlineSingle

# --- Apply GaNStyle to every "2022 Campaign Dates ..." run (4 occurrences) ---

$find1 = $d.Content
$find1.Find.ClearFormatting()
$find1.Find.Text = " 2022 Campaign Dates that use Hercules constellation: June 13-22, July 12-21, August 10-19"
while ($find1.Find.Execute()) {
    $find1.Style = "GaNStyle"
    $find1.Collapse(0)
}

# --- Apply GaNParagraph to the campaign description run (1 occurrence) ---

$find2 = $d.Content
$find2.Find.ClearFormatting()
$find2.Find.Text = "You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Hercules constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky."
if ($find2.Find.Execute()) {
    $find2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the amper.ped.muni.cz link run (1 occurrence) ---

$find3 = $d.Content
$find3.Find.ClearFormatting()
$find3.Find.Text = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
if ($find3.Find.Execute()) {
    $find3.Style = "GaNLinks"
}
